# Auto-generated edit script for toilet_coop_2022-07-23.xlsx
# Re-orders several scraped product rows to match a re-crawl at a later
# timestamp, and refreshes the 'timestamp' column (O) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell {
    # Forces text storage (NumberFormat @) so numeric-looking strings
    # like ids/prices ('6346771', '15.50') are not coerced to numbers.
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Set-NumCell {
    param($ws, $ref, $val)
    $ws.Range($ref).Value = $val
}

function Clear-Cell {
    param($ws, $ref)
    $ws.Range($ref).Value = $null
}

# Row 5 <- content previously shown at row 7
Set-TextCell $ws "A5" "6568452"
Set-TextCell $ws "B5" "Super Soft Premium Mandel feucht 4x  50ST"
Set-TextCell $ws "C5" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-premium-mandel-feucht/p/6568452"
Set-TextCell $ws "D5" "4x 50ST"
Set-TextCell $ws "G5" "Super Soft"
Set-TextCell $ws "H5" "6.65"
Set-TextCell $ws "I5" "0.03/1ST"
Set-TextCell $ws "J5" "Preis pro 1 Stück"
Set-TextCell $ws "K5" "0.03"
Set-TextCell $ws "L5" "1ST"
Set-TextCell $ws "M5" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
Set-TextCell $ws "N5" "Super Soft Premium Mandel feucht 4x  50ST 43% Aktion 6.65 Schweizer Franken statt 11.80 Schweizer Franken"
Set-NumCell $ws "E5" 8
Set-NumCell $ws "F5" 3.5

# Row 6 <- content previously shown at row 5
Set-TextCell $ws "A6" "6346771"
Set-TextCell $ws "B6" "Tempo Toilettenpapier Premium 4-lagig 16 Rollen"
Set-TextCell $ws "C6" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tempo-toilettenpapier-premium-4-lagig-16-rollen/p/6346771"
Set-TextCell $ws "D6" "16Rol"
Set-TextCell $ws "G6" "Tempo"
Set-TextCell $ws "H6" "15.50"
Set-TextCell $ws "I6" "0.97/1Rol"
Set-TextCell $ws "J6" "Preis pro 1 Rolle"
Set-TextCell $ws "K6" "0.97"
Set-TextCell $ws "L6" "1Rol"
Set-TextCell $ws "M6" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
Set-TextCell $ws "N6" "Tempo Toilettenpapier Premium 4-lagig 16 Rollen 15.50 Schweizer Franken"
Set-NumCell $ws "E6" 7
Set-NumCell $ws "F6" 4

# Row 7 <- content previously shown at row 6
Set-TextCell $ws "A7" "4947421"
Set-TextCell $ws "B7" "Oecoplan Taschentuch Calendula Box"
Set-TextCell $ws "C7" "/de/inspiration-geschenke/saisonale-promotionen/gesundheit/oecoplan-taschentuch-calendula-box/p/4947421"
Set-TextCell $ws "D7" "80ST"
Set-TextCell $ws "G7" "Coop"
Set-TextCell $ws "H7" "2.30"
Set-TextCell $ws "I7" "0.03/1ST"
Set-TextCell $ws "J7" "Preis pro 1 Stück"
Set-TextCell $ws "K7" "0.03"
Set-TextCell $ws "L7" "1ST"
Set-TextCell $ws "M7" "['inspiration-geschenke', 'saisonale-promotionen', 'gesundheit']"
Set-TextCell $ws "N7" "Oecoplan Taschentuch Calendula Box 2.30 Schweizer Franken"
Set-NumCell $ws "E7" 17
Set-NumCell $ws "F7" 4

# Row 10 <- content previously shown at row 12
Set-TextCell $ws "A10" "6691348"
Set-TextCell $ws "B10" "Super Soft Aloe Vera feucht FSC 4x  60ST"
Set-TextCell $ws "C10" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-aloe-vera-feucht-fsc/p/6691348"
Set-TextCell $ws "D10" "4x 60ST"
Set-TextCell $ws "G10" "Super Soft"
Set-TextCell $ws "H10" "6.65"
Set-TextCell $ws "I10" "0.03/1ST"
Set-TextCell $ws "J10" "Preis pro 1 Stück"
Set-TextCell $ws "K10" "0.03"
Set-TextCell $ws "L10" "1ST"
Set-TextCell $ws "M10" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
Set-TextCell $ws "N10" "Super Soft Aloe Vera feucht FSC 4x  60ST 43% Aktion 6.65 Schweizer Franken statt 11.80 Schweizer Franken"
Set-NumCell $ws "E10" 16
Set-NumCell $ws "F10" 1

# Row 11 <- content previously shown at row 10
Set-TextCell $ws "A11" "6695141"
Set-TextCell $ws "B11" "Prix Garantie feuchtes Toilettenpapier 2x70 Stück"
Set-TextCell $ws "C11" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/prix-garantie-feuchtes-toilettenpapier-2x70-stueck/p/6695141"
Set-TextCell $ws "D11" "140ST"
Set-TextCell $ws "G11" "Coop"
Set-TextCell $ws "H11" "2.40"
Set-TextCell $ws "I11" "0.02/1ST"
Set-TextCell $ws "J11" "Preis pro 1 Stück"
Set-TextCell $ws "K11" "0.02"
Set-TextCell $ws "L11" "1ST"
Set-TextCell $ws "M11" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
Set-TextCell $ws "N11" "Prix Garantie feuchtes Toilettenpapier 2x70 Stück 2.40 Schweizer Franken"
Set-NumCell $ws "E11" 5
Set-NumCell $ws "F11" 3.5

# Row 12 <- content previously shown at row 11
Set-TextCell $ws "A12" "6873015"
Set-TextCell $ws "B12" "Hakle Toilettenpapier Sagenhafte Sauberkeit 3-lagig 12 Rollen"
Set-TextCell $ws "C12" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/hakle-toilettenpapier-sagenhafte-sauberkeit-3-lagig-12-rollen/p/6873015"
Set-TextCell $ws "D12" "12Rol"
Set-TextCell $ws "G12" "Hakle"
Set-TextCell $ws "H12" "10.95"
Set-TextCell $ws "I12" "0.91/1Rol"
Set-TextCell $ws "J12" "Preis pro 1 Rolle"
Set-TextCell $ws "K12" "0.91"
Set-TextCell $ws "L12" "1Rol"
Set-TextCell $ws "M12" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
Set-TextCell $ws "N12" "Hakle Toilettenpapier Sagenhafte Sauberkeit 3-lagig 12 Rollen 10.95 Schweizer Franken"
Clear-Cell $ws "E12"
Set-NumCell $ws "F12" 0

# Row 16 <- content previously shown at row 17
Set-TextCell $ws "A16" "6346813"
Set-TextCell $ws "B16" "Tempo Toilettenpapier Premium 4-lagig 9 Rollen"
Set-TextCell $ws "C16" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tempo-toilettenpapier-premium-4-lagig-9-rollen/p/6346813"
Set-TextCell $ws "D16" "9Rol"
Set-TextCell $ws "G16" "Tempo"
Set-TextCell $ws "H16" "8.80"
Set-TextCell $ws "I16" "0.98/1Rol"
Set-TextCell $ws "J16" "Preis pro 1 Rolle"
Set-TextCell $ws "K16" "0.98"
Set-TextCell $ws "L16" "1Rol"
Set-TextCell $ws "M16" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
Set-TextCell $ws "N16" "Tempo Toilettenpapier Premium 4-lagig 9 Rollen 8.80 Schweizer Franken"
Set-NumCell $ws "E16" 6
Set-NumCell $ws "F16" 4.5

# Row 17 <- content previously shown at row 16
Set-TextCell $ws "A17" "6834305"
Set-TextCell $ws "B17" "Zewa Wisch&amp;Weg Haushaltspapier weiss 4 Rollen"
Set-TextCell $ws "C17" "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/zewa-wisch-weg-haushaltspapier-weiss-4-rollen/p/6834305"
Set-TextCell $ws "D17" "192BLT"
Set-TextCell $ws "G17" "Zewa"
Set-TextCell $ws "H17" "5.50"
Clear-Cell $ws "I17"
Clear-Cell $ws "J17"
Clear-Cell $ws "K17"
Clear-Cell $ws "L17"
Set-TextCell $ws "M17" "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
Set-TextCell $ws "N17" "Zewa Wisch&amp;Weg Haushaltspapier weiss 4 Rollen 5.50 Schweizer Franken"
Set-NumCell $ws "E17" 1
Set-NumCell $ws "F17" 5

# Row 25 <- content previously shown at row 26
Set-TextCell $ws "A25" "6498679"
Set-TextCell $ws "B25" "Subito Haushaltspapier decor 12 Rollen"
Set-TextCell $ws "C25" "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/subito-haushaltspapier-decor-12-rollen/p/6498679"
Set-TextCell $ws "D25" "600BLT"
Set-TextCell $ws "G25" "subito"
Set-TextCell $ws "H25" "10.95"
Clear-Cell $ws "I25"
Clear-Cell $ws "J25"
Clear-Cell $ws "K25"
Clear-Cell $ws "L25"
Set-TextCell $ws "M25" "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
Set-TextCell $ws "N25" "Subito Haushaltspapier decor 12 Rollen 33% Aktion 10.95 Schweizer Franken statt 16.50 Schweizer Franken"
Clear-Cell $ws "E25"
Set-NumCell $ws "F25" 0

# Row 26 <- content previously shown at row 25
Set-TextCell $ws "A26" "6149220"
Set-TextCell $ws "B26" "Naturaline Baby Feuchttücher 72 Stk."
Set-TextCell $ws "C26" "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/naturaline-baby-feuchttuecher-72-stk/p/6149220"
Set-TextCell $ws "D26" "72ST"
Set-TextCell $ws "G26" "Coop"
Set-TextCell $ws "H26" "3.95"
Set-TextCell $ws "I26" "0.05/1ST"
Set-TextCell $ws "J26" "Preis pro 1 Stück"
Set-TextCell $ws "K26" "0.05"
Set-TextCell $ws "L26" "1ST"
Set-TextCell $ws "M26" "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
Set-TextCell $ws "N26" "Naturaline Baby Feuchttücher 72 Stk. 3.95 Schweizer Franken"
Set-NumCell $ws "E26" 7
Set-NumCell $ws "F26" 4.5

# Refresh the crawl timestamp for every data row (2-37)
$newTimestamp = "2022-07-23 20:58:56"
for ($r = 2; $r -le 37; $r++) {
    Set-TextCell $ws ("O" + $r) $newTimestamp
}
